$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds plain decimal-looking price strings (e.g. "1.000", "0.5359").
# Excel normally auto-converts such text into real numbers when a value is
# assigned, which would change both the stored type and the displayed text
# (e.g. "1.000" -> 1). To preserve these as literal text (matching the
# original workbook, where every cell in D/E is a text string), the column is
# temporarily switched to the Text number format before the values are written,
# then the style is reset back to the default "Normal" style afterwards so the
# cells keep their original (unstyled) appearance.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.213.14'
$ws.Range('E2').Value = '  -0.02%  '

$ws.Range('D3').Value = '1.904.08'
$ws.Range('E3').Value = '  +0.41%  '

$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('E5').Value = '  -0.42%  '

$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.09%  '

$ws.Range('D7').Value = '0.5359'
$ws.Range('E7').Value = '  +3.13%  '

$ws.Range('D8').Value = '0.3807'
$ws.Range('E8').Value = '  +1.13%  '

$ws.Range('D9').Value = '0.07288'
$ws.Range('E9').Value = '  +0.07%  '

$ws.Range('D10').Value = '22.18'
$ws.Range('E10').Value = '  +4.66%  '

$ws.Range('D11').Value = '0.9051'
$ws.Range('E11').Value = '  +0.51%  '

$ws.Range('D12').Value = '0.08199'
$ws.Range('E12').Value = '  +0.43%  '

$ws.Range('D13').Value = '95.88'
$ws.Range('E13').Value = '  -0.68%  '

$ws.Range('D14').Value = '5.341'
$ws.Range('E14').Value = '  +1.19%  '

$ws.Range('E15').Value = '  -0.12%  '

$ws.Range('D16').Value = '14.84'
$ws.Range('E16').Value = '  +2.08%  '

$ws.Range('D17').Value = '0.000008662'
$ws.Range('E17').Value = '  +0.57%  '

$ws.Range('D18').Value = '1.000'
$ws.Range('E18').Value = '  -0.10%  '

$ws.Range('D19').Value = '27.238.10'
$ws.Range('E19').Value = '  -0.05%  '

$ws.Range('D20').Value = '5.043'
$ws.Range('E20').Value = '  -0.82%  '

$ws.Range('D21').Value = '1.082.11'
$ws.Range('E21').Value = '  -42.98%  '

$ws.Range('E22').Value = '  +0.84%  '

$ws.Range('D23').Value = '6.519'
$ws.Range('E23').Value = '  +1.91%  '

$ws.Range('D24').Value = '149.33'
$ws.Range('E24').Value = '  +1.38%  '

$ws.Range('D25').Value = '2.293'
$ws.Range('E25').Value = '  -0.19%  '

$ws.Range('D26').Value = '18.37'
$ws.Range('E26').Value = '  +1.02%  '

$ws.Range('D27').Value = '1.746'
$ws.Range('E27').Value = '  +0.10%  '

$ws.Range('D28').Value = '116.74'
$ws.Range('E28').Value = '  +1.36%  '

$ws.Range('D29').Value = '4.813'

$ws.Range('D30').Value = '4.724'
$ws.Range('E30').Value = '  -4.67%  '

$ws.Range('D31').Value = '0.09219'
$ws.Range('E31').Value = '  -0.13%  '

$ws.Range('D32').Value = '0.8293'
$ws.Range('E32').Value = '  +4.28%  '

$ws.Range('D33').Value = '0.05076'
$ws.Range('E33').Value = '  +0.87%  '

$ws.Range('E34').Value = '  -0.23%  '

$ws.Range('D35').Value = '3.005'
$ws.Range('E35').Value = '  +1.51%  '

$ws.Range('D36').Value = '3.336'
$ws.Range('E36').Value = '  -3.28%  '

$ws.Range('D37').Value = '2.672'
$ws.Range('E37').Value = '  +2.95%  '

$ws.Range('D38').Value = '0.5818'
$ws.Range('E38').Value = '  +2.56%  '

$ws.Range('D39').Value = '0.02004'
$ws.Range('E39').Value = '  +1.01%  '

$ws.Range('D40').Value = '1.076'
$ws.Range('E40').Value = '  +0.39%  '

$ws.Range('D41').Value = '9.320'
$ws.Range('E41').Value = '  +4.14%  '

$ws.Range('E42').Value = '  +0.98%  '

$ws.Range('D43').Value = '117.26'
$ws.Range('E43').Value = '  +1.50%  '

$ws.Range('D44').Value = '0.5054'
$ws.Range('E44').Value = '  +3.75%  '

$ws.Range('D45').Value = '0.1523'
$ws.Range('E45').Value = '  +0.53%  '

$ws.Range('D46').Value = '0.9998'
$ws.Range('E46').Value = '  -0.13%  '

$ws.Range('D47').Value = '10.08'
$ws.Range('E47').Value = '  -0.03%  '

$ws.Range('D48').Value = '1.640'
$ws.Range('E48').Value = '  +1.08%  '

$ws.Range('D49').Value = '38.36'
$ws.Range('E49').Value = '  +0.42%  '

$ws.Range('D50').Value = '0.06155'
$ws.Range('E50').Value = '  +3.62%  '

$ws.Range('D51').Value = '63.51'
$ws.Range('E51').Value = '  +0.11%  '

# Restore the default style on column D now that the values are safely stored as text.
$ws.Range("D2:D51").Style = "Normal"

Write-Host "Done"